$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1778656126482213
$ws.Range("C2").Value = 0.5968379446640316
$ws.Range("J2").Value = 0.01185770750988142
$ws.Range("O2").Value = 0.003952569169960474
$ws.Range("P2").Value = 0.1106719367588933
$ws.Range("S2").Value = 0.09881422924901186
$ws.Range("B3").Value = 0.0130718954248366
$ws.Range("C3").Value = 0.0261437908496732
$ws.Range("J3").Value = 0.0130718954248366
$ws.Range("P3").Value = 0.7712418300653595
$ws.Range("S3").Value = 0.1764705882352941
$ws.Range("J4").Value = 0.02380952380952381
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.3095238095238095
$ws.Range("B6").Value = 0.04504504504504504
$ws.Range("D6").Value = 0.009009009009009009
$ws.Range("F6").Value = 0.0990990990990991
$ws.Range("J6").Value = 0.2207207207207207
$ws.Range("O6").Value = 0.02702702702702703
$ws.Range("Q6").Value = 0.1891891891891892
$ws.Range("R6").Value = 0.05855855855855856
$ws.Range("S6").Value = 0.3513513513513514
$ws.Range("B7").Value = 0.08490566037735849
$ws.Range("D7").Value = 0.01886792452830189
$ws.Range("F7").Value = 0.03773584905660377
$ws.Range("J7").Value = 0.1179245283018868
$ws.Range("O7").Value = 0.02358490566037736
$ws.Range("Q7").Value = 0.1933962264150944
$ws.Range("R7").Value = 0.05660377358490566
$ws.Range("S7").Value = 0.4669811320754717
$ws.Range("B8").Value = 0.0851063829787234
$ws.Range("D8").Value = 0.02600472813238771
$ws.Range("E8").Value = 0.002364066193853428
$ws.Range("F8").Value = 0.06855791962174941
$ws.Range("J8").Value = 0.08037825059101655
$ws.Range("O8").Value = 0.01182033096926714
$ws.Range("Q8").Value = 0.2127659574468085
$ws.Range("R8").Value = 0.07092198581560284
$ws.Range("S8").Value = 0.442080378250591
$ws.Range("B9").Value = 0.09523809523809523
$ws.Range("D9").Value = 0.009523809523809525
$ws.Range("F9").Value = 0.0761904761904762
$ws.Range("J9").Value = 0.09523809523809523
$ws.Range("O9").Value = 0.0380952380952381
$ws.Range("Q9").Value = 0.1619047619047619
$ws.Range("R9").Value = 0.1095238095238095
$ws.Range("S9").Value = 0.4142857142857143
$ws.Range("B10").Value = 0.1191658391261172
$ws.Range("D10").Value = 0.02482621648460774
$ws.Range("F10").Value = 0.07845084409136048
$ws.Range("J10").Value = 0.0903674280039722
$ws.Range("O10").Value = 0.009930486593843098
$ws.Range("Q10").Value = 0.1996027805362463
$ws.Range("R10").Value = 0.08341608738828203
$ws.Range("S10").Value = 0.394240317775571
$ws.Range("G11").Value = 0.1623376623376623
$ws.Range("J11").Value = 0.06493506493506493
$ws.Range("K11").Value = 0.1948051948051948
$ws.Range("L11").Value = 0.5584415584415584
$ws.Range("S11").Value = 0.01948051948051948
$ws.Range("G12").Value = 0.7237569060773481
$ws.Range("J12").Value = 0.1878453038674033
$ws.Range("K12").Value = 0.01657458563535912
$ws.Range("L12").Value = 0.03867403314917127
$ws.Range("S12").Value = 0.03314917127071823
$ws.Range("G13").Value = 0.7777777777777778
$ws.Range("J13").Value = 0.2222222222222222
$ws.Range("F15").Value = 0.01408450704225352
$ws.Range("H15").Value = 0.1596244131455399
$ws.Range("I15").Value = 0.07981220657276995
$ws.Range("J15").Value = 0.352112676056338
$ws.Range("K15").Value = 0.0892018779342723
$ws.Range("M15").Value = 0.01408450704225352
$ws.Range("O15").Value = 0.09389671361502347
$ws.Range("S15").Value = 0.1971830985915493
$ws.Range("F16").Value = 0.01785714285714286
$ws.Range("H16").Value = 0.1845238095238095
$ws.Range("I16").Value = 0.07738095238095238
$ws.Range("J16").Value = 0.3571428571428572
$ws.Range("K16").Value = 0.119047619047619
$ws.Range("M16").Value = 0.04761904761904762
$ws.Range("O16").Value = 0.07738095238095238
$ws.Range("S16").Value = 0.119047619047619
$ws.Range("F17").Value = 0.02457002457002457
$ws.Range("H17").Value = 0.1793611793611794
$ws.Range("I17").Value = 0.1105651105651106
$ws.Range("J17").Value = 0.3857493857493858
$ws.Range("K17").Value = 0.09828009828009827
$ws.Range("M17").Value = 0.01228501228501228
$ws.Range("O17").Value = 0.05896805896805897
$ws.Range("S17").Value = 0.1302211302211302
$ws.Range("F18").Value = 0.0124223602484472
$ws.Range("H18").Value = 0.1925465838509317
$ws.Range("I18").Value = 0.124223602484472
$ws.Range("J18").Value = 0.3354037267080746
$ws.Range("K18").Value = 0.08074534161490683
$ws.Range("M18").Value = 0.01863354037267081
$ws.Range("N18").Value = 0.0124223602484472
$ws.Range("O18").Value = 0.08074534161490683
$ws.Range("S18").Value = 0.1428571428571428
$ws.Range("F19").Value = 0.01745635910224439
$ws.Range("H19").Value = 0.2128013300083126
$ws.Range("I19").Value = 0.09642560266001662
$ws.Range("J19").Value = 0.3183707398171239
$ws.Range("K19").Value = 0.1230257689110557
$ws.Range("M19").Value = 0.02161263507896924
$ws.Range("N19").Value = 0.0008312551953449709
$ws.Range("O19").Value = 0.07315045719035744
$ws.Range("S19").Value = 0.1363258520365752
